$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.102.37"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.937.73"
$ws.Range("E3").Value = "  +3.49%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.18"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5083"
$ws.Range("E7").Value = "  +0.90%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4043"
$ws.Range("E8").Value = "  +2.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08410"
$ws.Range("E9").Value = "  +2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  +2.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.27"
$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.02"
$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.423"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").Value = "1.924.89"
$ws.Range("E14").Value = "  +3.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.294"
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.01"
$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001099"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06509"
$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.54"
$ws.Range("E20").Value = "  +2.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.976"
$ws.Range("E22").Value = "  +2.34%  "

$ws.Range("D23").Value = "30.120.76"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.190"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("D26").Value = "2.143.70"
$ws.Range("E26").Value = "  +2.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.02"
$ws.Range("E27").Value = "  +2.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.60"
$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.276"
$ws.Range("E29").Value = "  +2.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.87"
$ws.Range("E30").Value = "  +2.11%  "

$ws.Range("E31").Value = "  +6.27%  "

$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.991"
$ws.Range("E33").Value = "  +0.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.774"
$ws.Range("E34").Value = "  +2.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02457"
$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.326"
$ws.Range("E36").Value = "  +1.81%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.261"
$ws.Range("E37").Value = "  +7.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06471"
$ws.Range("E38").Value = "  +1.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2155"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6494"
$ws.Range("E40").Value = "  +2.99%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.720"
$ws.Range("E41").Value = "  +2.60%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.72"
$ws.Range("E42").Value = "  +3.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.221"
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6078"
$ws.Range("E44").Value = "  +2.78%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.33"
$ws.Range("E45").Value = "  +2.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.173"
$ws.Range("E46").Value = "  +3.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.629"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.60"
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.212"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.135"
$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.22"
$ws.Range("E51").Value = "  +1.08%  "
